$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '67.638.92'
$ws.Cells.Item(2, 5).Value = '  -2.38%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.272.05'
$ws.Cells.Item(3, 5).Value = '  -4.90%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.09%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '595.29'
$ws.Cells.Item(5, 5).Value = '  -2.69%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '149.24'
$ws.Cells.Item(6, 5).Value = '  -10.89%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.262.60'
$ws.Cells.Item(8, 5).Value = '  -4.93%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.543'
$ws.Cells.Item(9, 5).Value = '  -8.42%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -11.62%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '6.63'
$ws.Cells.Item(11, 5).Value = '  -6.11%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.504'
$ws.Cells.Item(12, 5).Value = '  -10.74%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).Value = '0.0000246'
$ws.Cells.Item(13, 5).Value = '  -8.73%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).Value = '38.34'
$ws.Cells.Item(14, 5).Value = '  -13.27%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.793.65'
$ws.Cells.Item(15, 5).Value = '  -5.05%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '67.651.94'
$ws.Cells.Item(16, 5).Value = '  -2.45%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.277.11'
$ws.Cells.Item(17, 5).Value = '  -4.89%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'TRON'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(18, 4).Value = '0.114'
$ws.Cells.Item(18, 5).Value = '  -5.81%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'BitcoinCash'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(19, 4).Value = '528.85'
$ws.Cells.Item(19, 5).Value = '  -9.01%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '7.12'
$ws.Cells.Item(20, 5).Value = '  -12.57%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '14.92'
$ws.Cells.Item(21, 5).Value = '  -12.81%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.753'
$ws.Cells.Item(22, 5).Value = '  -10.97%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '7.81'
$ws.Cells.Item(23, 5).Value = '  -11.75%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '85.64'
$ws.Cells.Item(24, 5).Value = '  -10.54%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '13.46'
$ws.Cells.Item(25, 5).Value = '  -11.09%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.13%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '3.23'
$ws.Cells.Item(27, 5).Value = '  -10.91%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'ImmutableX'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(28, 4).Value = '2.14'
$ws.Cells.Item(28, 5).Value = '  -11.45%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'RenderToken'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(29, 4).Value = '8.02'
$ws.Cells.Item(29, 5).Value = '  -6.76%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '28.98'
$ws.Cells.Item(30, 5).Value = '  -11.35%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.20'
$ws.Cells.Item(31, 5).Value = '  -3.16%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '2.65'
$ws.Cells.Item(32, 5).Value = '  -4.43%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '6.61'
$ws.Cells.Item(33, 5).Value = '  -15.31%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.70'
$ws.Cells.Item(34, 5).Value = '  -12.78%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.02%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'OKB'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(36, 4).Value = '56.56'
$ws.Cells.Item(36, 5).Value = '  +1.52%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '513.07'
$ws.Cells.Item(37, 5).Value = '  -10.78%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.0444'
$ws.Cells.Item(38, 5).Value = '  -6.43%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.0854'
$ws.Cells.Item(39, 5).Value = '  -10.65%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '8.92'
$ws.Cells.Item(40, 5).Value = '  -15.07%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -10.90%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '2.76'
$ws.Cells.Item(42, 5).Value = '  -12.26%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '2.930.84'
$ws.Cells.Item(43, 5).Value = '  -9.49%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '0.267'
$ws.Cells.Item(44, 5).Value = '  -9.57%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '0.0₃0579'
$ws.Cells.Item(45, 5).Value = '  -15.19%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.20'
$ws.Cells.Item(46, 5).Value = '  -8.00%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '26.56'
$ws.Cells.Item(47, 5).Value = '  -14.61%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.08%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.31'
$ws.Cells.Item(49, 5).Value = '  -16.27%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -10.20%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '123.77'
$ws.Cells.Item(51, 5).Value = '  -7.54%  '
